$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Add the new 2025 row of data just below the existing table (row 14)
$ws.Range("A14").Value = 2025
$ws.Range("B14").Value = 30
$ws.Range("C14").Value = 386

# Update selection on the Graph sheet to match the recorded edit
$graph = $wb.Worksheets.Item("Graph")
$graph.Range("J22").Select()
